$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '55.878.60'
$ws.Range('E2').Value = '  +2.62%  '
$ws.Range('D3').Value = '2.485.90'
$ws.Range('E3').Value = '  +2.79%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'483.06"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.79%  '
$ws.Range('D6').Value = "'144.09"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +10.54%  '
$ws.Range('E7').Value = '  +0.32%  '
$ws.Range('E8').Value = '  +3.62%  '
$ws.Range('D9').Value = '2.505.77'
$ws.Range('E9').Value = '  +3.17%  '
$ws.Range('E10').Value = '  +4.85%  '
$ws.Range('D11').Value = "'0.0967"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.15%  '
$ws.Range('E12').Value = '  +3.48%  '
$ws.Range('E13').Value = '  +0.97%  '
$ws.Range('D14').Value = '2.916.18'
$ws.Range('E14').Value = '  +3.25%  '
$ws.Range('D15').Value = '55.760.17'
$ws.Range('E15').Value = '  +2.83%  '
$ws.Range('D16').Value = "'20.85"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.86%  '
$ws.Range('E17').Value = '  +3.84%  '
$ws.Range('D18').Value = '2.501.15'
$ws.Range('E18').Value = '  +3.56%  '
$ws.Range('D19').Value = "'4.42"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.59%  '
$ws.Range('D20').Value = "'10.21"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +8.79%  '
$ws.Range('D21').Value = "'317.38"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.73%  '
$ws.Range('D22').Value = "'0.999"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('D23').Value = "'5.77"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +6.58%  '
$ws.Range('D24').Value = "'58.13"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.56%  '
$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').Value = "'0.408"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.90%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').Value = "'0.165"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.80%  '
$ws.Range('D27').Value = "'0.998"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.65%  '
$ws.Range('D28').Value = '2.606.36'
$ws.Range('E28').Value = '  +4.44%  '
$ws.Range('E29').Value = '  +3.74%  '
$ws.Range('D30').Value = '0.0₃0777'
$ws.Range('E30').Value = '  +9.43%  '
$ws.Range('E31').Value = '  +0.37%  '
$ws.Range('D32').Value = "'148.19"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.37%  '
$ws.Range('E33').Value = '  +3.25%  '
$ws.Range('E34').Value = '  +5.67%  '
$ws.Range('E35').Value = '  +3.00%  '
$ws.Range('E36').Value = '  +8.00%  '
$ws.Range('E37').Value = '  +3.58%  '
$ws.Range('D38').Value = "'0.853"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.13%  '
$ws.Range('E39').Value = '  +1.40%  '
$ws.Range('D40').Value = "'3.51"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.33%  '
$ws.Range('D41').Value = "'0.997"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.51%  '
$ws.Range('E42').Value = '  +1.73%  '
$ws.Range('D43').Value = "'0.0550"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.58%  '
$ws.Range('D44').Value = "'1.31"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.61%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = "'4.75"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +10.26%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').Value = "'259.07"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +18.76%  '
$ws.Range('D47').Value = "'10.16"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.14%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = "'0.0225"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.56%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = "'0.0900"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.67%  '
$ws.Range('D50').Value = '1.913.90'
$ws.Range('E50').Value = '  -2.89%  '
$ws.Range('D51').Value = "'17.56"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.32%  '
